$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A80").Value = "Southern Institute of Ecology, Vietnam Academy of Science and Technology"

$ws.Range("F1").Select()
$excel.ActiveWindow.ScrollRow = 67
